# Slide 5, shape "TextBox 6" (Shapes.Item(8)) contains the misspelled
# heading "Local Gornment Area Offences base on Police Region" split
# across three runs (the middle one flagged with a spell-check err="1").
# The target fixes the typo ("Gornment" -> "Government") and merges the
# three runs into a single run that keeps the first run's formatting.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(8)
$tr = $sh.TextFrame.TextRange

# Remove everything after "Local " (i.e. the 2nd and 3rd runs), leaving
# a single remaining run that carries the first run's formatting
# (sz=3200, dirty=0, accent2 fill) plus the paragraph's endParaRPr.
$len = $tr.Length
$tail = $tr.Characters(7, $len - 6)
$tail.Text = ""

# Re-set the full corrected text; because only one run now exists, the
# engine extends that single run rather than re-splitting into multiple
# runs, so the corrected text ends up in one run with the original
# formatting and no err="1" spell-flag.
$tr.Text = "Local Government Area Offences base on Police Region"
